$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($sheet, $col, $row, $text)
    $cell = $sheet.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws 4 2 "62.622.70"
Set-TextValue $ws 5 2 "  +1.13%  "
Set-TextValue $ws 4 3 "3.017.25"
Set-TextValue $ws 5 3 "  -0.41%  "
Set-TextValue $ws 5 4 "  -0.20%  "
Set-TextValue $ws 4 5 "546.78"
Set-TextValue $ws 5 5 "  -0.98%  "
Set-TextValue $ws 4 6 "140.18"
Set-TextValue $ws 5 6 "  +3.11%  "
Set-TextValue $ws 4 7 "0.998"
Set-TextValue $ws 5 7 "  -0.19%  "
Set-TextValue $ws 4 8 "3.016.03"
Set-TextValue $ws 5 8 "  -0.29%  "
Set-TextValue $ws 4 9 "0.492"
Set-TextValue $ws 5 9 "  -1.18%  "
Set-TextValue $ws 4 10 "6.99"
Set-TextValue $ws 5 10 "  +14.90%  "
Set-TextValue $ws 4 11 "0.149"
Set-TextValue $ws 5 11 "  -0.02%  "
Set-TextValue $ws 4 12 "0.449"
Set-TextValue $ws 5 12 "  -0.66%  "
Set-TextValue $ws 4 13 "0.0000222"
Set-TextValue $ws 5 13 "  -0.55%  "
Set-TextValue $ws 4 14 "34.27"
Set-TextValue $ws 5 14 "  -0.63%  "
Set-TextValue $ws 4 15 "3.475.48"
Set-TextValue $ws 5 15 "  -1.22%  "
Set-TextValue $ws 4 16 "62.544.73"
Set-TextValue $ws 5 16 "  +0.78%  "
Set-TextValue $ws 4 17 "3.009.73"
Set-TextValue $ws 5 17 "  -0.65%  "
Set-TextValue $ws 5 18 "  -1.97%  "
Set-TextValue $ws 4 19 "6.61"
Set-TextValue $ws 5 19 "  -1.12%  "
Set-TextValue $ws 4 20 "474.07"
Set-TextValue $ws 5 20 "  -0.24%  "
Set-TextValue $ws 4 21 "13.50"
Set-TextValue $ws 5 21 "  +1.40%  "
Set-TextValue $ws 4 22 "0.658"
Set-TextValue $ws 5 22 "  -2.87%  "
Set-TextValue $ws 4 23 "7.22"
Set-TextValue $ws 5 23 "  +1.48%  "
Set-TextValue $ws 4 24 "79.80"
Set-TextValue $ws 5 24 "  -0.79%  "
Set-TextValue $ws 4 25 "12.70"
Set-TextValue $ws 5 25 "  +4.29%  "
Set-TextValue $ws 4 26 "1.00"
Set-TextValue $ws 5 26 "  +0.26%  "
Set-TextValue $ws 4 27 "2.74"
Set-TextValue $ws 5 27 "  -0.16%  "
Set-TextValue $ws 4 28 "7.70"
Set-TextValue $ws 5 28 "  -1.60%  "
Set-TextValue $ws 4 29 "2.02"
Set-TextValue $ws 5 29 "  +5.56%  "
Set-TextValue $ws 4 30 "0.996"
Set-TextValue $ws 5 30 "  -0.43%  "
Set-TextValue $ws 4 31 "25.62"
Set-TextValue $ws 5 31 "  -0.85%  "
Set-TextValue $ws 4 32 "1.13"
Set-TextValue $ws 5 32 "  -2.44%  "
Set-TextValue $ws 4 33 "2.37"
Set-TextValue $ws 5 33 "  +2.11%  "
Set-TextValue $ws 4 34 "5.61"
Set-TextValue $ws 5 34 "  +2.27%  "
Set-TextValue $ws 4 35 "54.75"
Set-TextValue $ws 5 35 "  -1.60%  "
Set-TextValue $ws 4 36 "5.88"
Set-TextValue $ws 5 36 "  -0.86%  "
Set-TextValue $ws 4 37 "455.40"
Set-TextValue $ws 5 37 "  -1.07%  "
Set-TextValue $ws 4 38 "0.0816"
Set-TextValue $ws 5 38 "  +1.88%  "
Set-TextValue $ws 4 39 "0.0395"
Set-TextValue $ws 5 39 "  +2.53%  "
Set-TextValue $ws 4 40 "2.969.35"
Set-TextValue $ws 5 40 "  -7.85%  "
Set-TextValue $ws 4 41 "0.115"
Set-TextValue $ws 5 41 "  -3.23%  "
Set-TextValue $ws 4 42 "8.11"
Set-TextValue $ws 5 42 "  -0.78%  "
Set-TextValue $ws 4 43 "2.58"
Set-TextValue $ws 5 43 "  +4.62%  "
Set-TextValue $ws 4 44 "27.03"
Set-TextValue $ws 5 44 "  +3.74%  "
Set-TextValue $ws 4 47 "2.02"
Set-TextValue $ws 5 47 "  +0.90%  "
Set-TextValue $ws 4 48 "0.110"
Set-TextValue $ws 5 48 "  +0.84%  "
Set-TextValue $ws 4 49 "115.50"
Set-TextValue $ws 5 49 "  -2.60%  "
Set-TextValue $ws 4 50 "0.0₃0498"
Set-TextValue $ws 5 50 "  +0.08%  "
Set-TextValue $ws 4 51 "2.03"
Set-TextValue $ws 5 51 "  +0.14%  "
Set-TextValue $ws 2 45 "TheGraph"
Set-TextValue $ws 3 45 "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws 4 45 "0.250"
Set-TextValue $ws 5 45 "  +1.92%  "
Set-TextValue $ws 2 46 "USDe"
Set-TextValue $ws 3 46 "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws 4 46 "1.00"
Set-TextValue $ws 5 46 "  -0.04%  "
